# Generate Report for Handback
# - Row 7 (749453a4-673a-4fca-bc64-2508a0f056ab) status flips from
#   "Ready for handoff" to "Handback transform failed" on the Overview
#   sheet and on both the zh-cn and de-de detail sheets.
# - Each detail sheet records an Error Detail message explaining the
#   handback/handoff file-name mismatch for that locale.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

$newStatus = "Handback transform failed"

$wsOverview.Range("B7").Value = $newStatus
$wsOverview.Range("C7").Value = $newStatus

$wsZhCn.Range("C7").Value = $newStatus
$wsDeDe.Range("C7").Value = $newStatus

$wsZhCn.Range("L7").Value = "Handback file name: eohk1115.kwf is different with handoff file name: 749453a4-673a-4fca-bc64-2508a0f056ab.93d288b2f67fa21e39ce08701c1f746ac7de8d6e.zh-cn."
$wsDeDe.Range("L7").Value = "Handback file name: eohk1115.kwf is different with handoff file name: 749453a4-673a-4fca-bc64-2508a0f056ab.93d288b2f67fa21e39ce08701c1f746ac7de8d6e.de-de."
